# completed synopsis and results after receiving final data
#
# - Column A (the running reference-list number) is shifted up by one for
#   every study row from "Bullock et al. 2021" (row 8) through
#   "Ziegler M. & Rondot P. (1999)" (row 66): 9->10, 11->12, 12->13, ...,
#   117->118. Rows 2-7 (numbers 1,2,3,5,6,8) are untouched.
# - Mizuno et al. 1995 (row 36) loses its extra pergolide-levodopa
#   baseline/follow-up outcome columns (F36, G36) - they are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 8; $r -le 66; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = [int]($cell.Value())
    $new = $old + 1

    # Force text storage (the column holds numbers-as-text, t="s" in the
    # xlsx, not numeric cells) and then drop back to the default "Normal"
    # style so no stray per-cell formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = [string]$new
    $cell.Style = "Normal"
}

# Clear the now-removed Mizuno "pergolide-levodopa" baseline/follow-up cells.
$clearRange = $ws.Range("F36:G36")
$clearRange.Font.Bold = $false
$clearRange.Value = ""
$clearRange.Style = "Normal"
